# Updated symbol list on Tue Jan 17 20:56:30 UTC 2023 with GitHub Actions
#
# The upstream scraper refreshed the "Price" (column D) and "Volume(1h)"
# (column E) columns for a number of rows in the crypto price sheet.
# All of these cells were authored as literal text (inlineStr) in the
# original workbook, e.g. "302.19" and "0.55%" are text, not a number or
# a percentage. We therefore write each value with a leading apostrophe
# so Excel stores it as text instead of auto-converting it to a number /
# percentage, and immediately clear the resulting cell formatting so the
# "quote prefix" text style doesn't stick around on the cell (keeping
# every cell's style identical to before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "302.19" },
    @{ Cell = "E2"; Value = "0.55%" },
    @{ Cell = "E3"; Value = "1.07%" },
    @{ Cell = "D4"; Value = "5.003" },
    @{ Cell = "E4"; Value = "-2.76%" },
    @{ Cell = "D5"; Value = "0.07911" },
    @{ Cell = "E5"; Value = "-3.04%" },
    @{ Cell = "D6"; Value = "2.120" },
    @{ Cell = "E6"; Value = "-16.03%" },
    @{ Cell = "D7"; Value = "7.870" },
    @{ Cell = "D8"; Value = "3.802" },
    @{ Cell = "E8"; Value = "-1.68%" },
    @{ Cell = "D9"; Value = "0.9285" },
    @{ Cell = "E9"; Value = "0.36%" },
    @{ Cell = "D10"; Value = "0.1751" },
    @{ Cell = "E10"; Value = "-0.65%" },
    @{ Cell = "D11"; Value = "0.07889" },
    @{ Cell = "E11"; Value = "6.41%" },
    @{ Cell = "D12"; Value = "0.08765" },
    @{ Cell = "E12"; Value = "-1.44%" },
    @{ Cell = "D13"; Value = "0.03142" },
    @{ Cell = "E13"; Value = "3.88%" },
    @{ Cell = "D14"; Value = "0.1003" },
    @{ Cell = "E14"; Value = "0.12%" },
    @{ Cell = "D15"; Value = "0.001532" },
    @{ Cell = "E15"; Value = "0.74%" },
    @{ Cell = "D16"; Value = "0.006022" },
    @{ Cell = "E16"; Value = "0.41%" },
    @{ Cell = "D17"; Value = "3.468" },
    @{ Cell = "E17"; Value = "-3.82%" },
    @{ Cell = "E18"; Value = "-0.31%" },
    @{ Cell = "D20"; Value = "0.1292" },
    @{ Cell = "E20"; Value = "-3.54%" },
    @{ Cell = "D21"; Value = "4.159" },
    @{ Cell = "E21"; Value = "2.03%" },
    @{ Cell = "D22"; Value = "0.1794" },
    @{ Cell = "E22"; Value = "6.83%" },
    @{ Cell = "D23"; Value = "0.04615" },
    @{ Cell = "E23"; Value = "-0.49%" },
    @{ Cell = "E24"; Value = "-0.61%" },
    @{ Cell = "D25"; Value = "0.004475" },
    @{ Cell = "E25"; Value = "-1.51%" },
    @{ Cell = "D26"; Value = "0.0001252" },
    @{ Cell = "E26"; Value = "4.48%" },
    @{ Cell = "D39"; Value = "0.01729" },
    @{ Cell = "E39"; Value = "-2.22%" },
    @{ Cell = "D40"; Value = "0.04822" },
    @{ Cell = "E40"; Value = "4.59%" },
    @{ Cell = "D41"; Value = "0.007351" },
    @{ Cell = "E41"; Value = "7.28%" },
    @{ Cell = "D42"; Value = "0.1364" },
    @{ Cell = "E42"; Value = "-1.14%" },
    @{ Cell = "D43"; Value = "0.002345" },
    @{ Cell = "E43"; Value = "6.20%" },
    @{ Cell = "D44"; Value = "0.01113" },
    @{ Cell = "E44"; Value = "12.99%" },
    @{ Cell = "D45"; Value = "0.00006037" },
    @{ Cell = "E45"; Value = "-2.10%" },
    @{ Cell = "E46"; Value = "0.14%" },
    @{ Cell = "D47"; Value = "0.003394" },
    @{ Cell = "E47"; Value = "-59.59%" },
    @{ Cell = "D48"; Value = "0.8205" },
    @{ Cell = "E48"; Value = "2.00%" },
    @{ Cell = "E49"; Value = "0.14%" },
    @{ Cell = "E50"; Value = "0.14%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.ClearFormats()
}
